$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update group name values (ID 0001/0002 prefix moved to suffix)
$ws.Range("B2").Value = "Ren_GP_0001"
$ws.Range("B3").Value = "Ren_GP_0002"

# Row 3 instructions / post-submission instructions now match row 2's text
$ws.Range("E3").Value = "This message is instructions text."
$ws.Range("F3").Value = "This message is post submission instructions text."

# Move active selection
[void]$ws.Range("B4").Select()
